# Auto-generated edit script applying scheduled-runner price/profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 210.97058
$ws.Range("I55").Value = 540.7
$ws.Range("J55").Value = 73.583336
$ws.Range("K55").Value = 540.7
$ws.Range("L55").Value = 73.583336
$ws.Range("M55").Value = -326.7
$ws.Range("N55").Value = -501.583336

$ws.Range("H58").Value = 3834.1667
$ws.Range("I58").Value = 751.25
$ws.Range("K58").Value = 2253.75
$ws.Range("M58").Value = -2103.75

$ws.Range("H80").Value = 10153065
$ws.Range("I80").Value = 857
$ws.Range("J80").Value = 16244390
$ws.Range("K80").Value = 2571
$ws.Range("L80").Value = 48733170
$ws.Range("M80").Value = -1573
$ws.Range("N80").Value = -48735166

$ws.Range("H83").Value = 10153065
$ws.Range("I83").Value = 857
$ws.Range("J83").Value = 16244390
$ws.Range("K83").Value = 7713
$ws.Range("L83").Value = 146199510
$ws.Range("M83").Value = -2721
$ws.Range("N83").Value = -146209494

$ws.Range("H100").Value = 52633500
$ws.Range("I100").Value = 166667500
$ws.Range("J100").Value = 2424.5386
$ws.Range("K100").Value = 166667500
$ws.Range("L100").Value = 2424.5386
$ws.Range("M100").Value = -166666959
$ws.Range("N100").Value = -3506.5386

$ws.Range("H103").Value = 625091.25
$ws.Range("I103").Value = 1250092.5
$ws.Range("J103").Value = 90
$ws.Range("K103").Value = 3750277.5
$ws.Range("L103").Value = 270
$ws.Range("M103").Value = -3749691.5
$ws.Range("N103").Value = -1442

$ws.Range("H113").Value = 37041704
$ws.Range("I113").Value = 76926850
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 76926850
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = -76923596
$ws.Range("N113").Value = -12008

$ws.Range("H116").Value = 4438.5625
$ws.Range("I116").Value = 1875.625
$ws.Range("J116").Value = 7001.5
$ws.Range("K116").Value = 1875.625
$ws.Range("L116").Value = 7001.5
$ws.Range("M116").Value = 1566.375
$ws.Range("N116").Value = -13885.5

$ws.Range("H118").Value = 895
$ws.Range("I118").Value = 895
$ws.Range("K118").Value = 2685
$ws.Range("M118").Value = -1028

$ws.Range("H132").Value = 4427.15
$ws.Range("I132").Value = 4574.278
$ws.Range("J132").Value = 3103
$ws.Range("K132").Value = 13722.834
$ws.Range("L132").Value = 9309
$ws.Range("M132").Value = -11192.834
$ws.Range("N132").Value = -14369

$ws.Range("H135").Value = 12198571
$ws.Range("I135").Value = 625.55884
$ws.Range("J135").Value = 71445736
$ws.Range("K135").Value = 5630.02956
$ws.Range("L135").Value = 643011624
$ws.Range("M135").Value = -3095.02956
$ws.Range("N135").Value = -643016694

$ws.Range("H137").Value = 37057.035
$ws.Range("I137").Value = 1044.2222
$ws.Range("K137").Value = 3132.6666
$ws.Range("M137").Value = -582.6665999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2540.11
$ws.Range("I32").Value = 2253.5227
$ws.Range("K32").Value = 2253.5227
$ws.Range("M32").Value = -1966.5227

$ws.Range("H61").Value = 1060470.5
$ws.Range("I61").Value = 1201733.6
$ws.Range("J61").Value = 997
$ws.Range("K61").Value = 1201733.6
$ws.Range("L61").Value = 997
$ws.Range("M61").Value = -1201521.6
$ws.Range("N61").Value = -1421

$ws.Range("H74").Value = 24391838
$ws.Range("I74").Value = 27779508
$ws.Range("K74").Value = 27779508
$ws.Range("M74").Value = -27778634

$ws.Range("H77").Value = 24391838
$ws.Range("I77").Value = 27779508
$ws.Range("K77").Value = 138897540
$ws.Range("M77").Value = -138893172

$ws.Range("H125").Value = 33851
$ws.Range("J125").Value = 33851
$ws.Range("L125").Value = 33851
$ws.Range("N125").Value = -43691

$ws.Range("H132").Value = 24193.912
$ws.Range("I132").Value = 2450.5264
$ws.Range("J132").Value = 127475
$ws.Range("K132").Value = 7351.5792
$ws.Range("L132").Value = 382425
$ws.Range("M132").Value = -4821.5792
$ws.Range("N132").Value = -387485

$ws.Range("H136").Value = 1060470.5
$ws.Range("I136").Value = 1201733.6
$ws.Range("J136").Value = 997
$ws.Range("K136").Value = 3605200.8
$ws.Range("L136").Value = 2991
$ws.Range("M136").Value = -3602650.8
$ws.Range("N136").Value = -8091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1847.7778
$ws.Range("I99").Value = 2310
$ws.Range("J99").Value = 1616.6666
$ws.Range("K99").Value = 2310
$ws.Range("L99").Value = 1616.6666
$ws.Range("M99").Value = -812
$ws.Range("N99").Value = -4612.6666

$ws.Range("H124").Value = 42000
$ws.Range("J124").Value = 42000
$ws.Range("L124").Value = 42000
$ws.Range("N124").Value = -51820

$ws.Range("H134").Value = 25182.445
$ws.Range("I134").Value = 26695.477
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 80086.431
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -77551.431
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 21742666
$ws.Range("I99").Value = 3269.2307
$ws.Range("K99").Value = 3269.2307
$ws.Range("M99").Value = -1771.2307

$ws.Range("H124").Value = 9662.5
$ws.Range("I124").Value = 9662.5
$ws.Range("K124").Value = 9662.5
$ws.Range("M124").Value = -7207.5

$ws.Range("H126").Value = 21742666
$ws.Range("I126").Value = 3269.2307
$ws.Range("K126").Value = 9807.6921
$ws.Range("M126").Value = -7337.6921

$ws.Range("H132").Value = 2474.3
$ws.Range("I132").Value = 1835
$ws.Range("J132").Value = 21014
$ws.Range("K132").Value = 5505
$ws.Range("L132").Value = 63042
$ws.Range("M132").Value = -2975
$ws.Range("N132").Value = -68102

$ws.Range("H134").Value = 844.2414
$ws.Range("I134").Value = 844.2414
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2532.7242
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 2.27579999999989
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 692.12
$ws.Range("I131").Value = 480
$ws.Range("J131").Value = 713.0989
$ws.Range("K131").Value = 1440
$ws.Range("L131").Value = 2139.2967
$ws.Range("M131").Value = 3600
$ws.Range("N131").Value = -12219.2967

$ws.Range("H141").Value = 2205.875
$ws.Range("I141").Value = 1807.0714
$ws.Range("K141").Value = 5421.2142
$ws.Range("M141").Value = -241.2142000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 6410484
$ws.Range("I107").Value = 303.55554
$ws.Range("J107").Value = 25641024
$ws.Range("K107").Value = 303.55554
$ws.Range("L107").Value = 25641024
$ws.Range("M107").Value = 1616.44446
$ws.Range("N107").Value = -25644864

$ws.Range("H132").Value = 87754.05499999999
$ws.Range("I132").Value = 63504.293
$ws.Range("K132").Value = 190512.879
$ws.Range("M132").Value = -187982.879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5966.5835
$ws.Range("I61").Value = 2825
$ws.Range("K61").Value = 2825
$ws.Range("M61").Value = -2623

$ws.Range("H113").Value = 5966.5835
$ws.Range("I113").Value = 2825
$ws.Range("K113").Value = 2825
$ws.Range("M113").Value = -655

$ws.Range("H122").Value = 819694.3
$ws.Range("I122").Value = 1636280.4
$ws.Range("K122").Value = 4908841.199999999
$ws.Range("M122").Value = -4906391.199999999

$ws.Range("H132").Value = 575173.1
$ws.Range("I132").Value = 603861.75
$ws.Range("K132").Value = 1811585.25
$ws.Range("M132").Value = -1809055.25
